# Update the dSF column (column F) values on the active worksheet.
# These are raw data re-pulls / recalculated means from the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -5
    4  = 6
    6  = -5
    9  = 9
    12 = -10
    14 = 0
    15 = -1
    16 = -1
    17 = -4
    23 = -1
    27 = 3
    31 = 5
    33 = 2
    34 = 1
    35 = -3
    36 = -4
    38 = -6
    39 = 0
    40 = 1
    43 = -1
    44 = 3
    45 = -5
    46 = 4
    47 = -1
    48 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
